$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$wsOther = $wb.Worksheets.Item("WAT09")

# --- Fill the new rows with the same "blank" formatting as row 2 (A:E) ---
$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A3:E19").PasteSpecial(-4122) | Out-Null
$ws.Range("A20:E21").PasteSpecial(-4122) | Out-Null
$ws.Range("A22:E32").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 20: Publons login test using Facebook ---
$ws.Range("A20").Value = "PUBLONS020"
$ws.Range("B20").Value = "OPQA-5885||OPQA-5888"
$ws.Range("C20").Value = "User Sign In  Publon's using valid facebook username and password from login page then it should redirected to Publon's Home page."
$ws.Range("D20").Value = "Y"

# --- Row 21: Publons login test using Gmail ---
$ws.Range("A21").Value = "PUBLONS021"
$ws.Range("B21").Value = "OPQA-5884||OPQA-5887"
$ws.Range("D21").Value = "Y"

# C21 description was pasted in from elsewhere (e.g. browser/Word), carrying its
# own rich-text formatting: Calibri 11, font color #212121, white fill,
# left/top aligned, vertical centered, wrapped text, with the sheet's
# standard thin border.
$wsOther.Range("A2").Copy() | Out-Null
$c21 = $ws.Range("C21")
$c21.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$c21.Value = "User Sign In Publon's using valid gmail username and password from login page then it should redirected to Publon's Home page."
$c21.Font.Color = 33 + 33 * 256 + 33 * 65536
$c21.Interior.Color = 255 + 255 * 256 + 255 * 65536
$c21.HorizontalAlignment = -4131
$c21.VerticalAlignment = -4108
$c21.WrapText = $true

# --- Sheet view / selection ---
$ws.Range("B12").Select() | Out-Null

Write-Host "done"
